$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.31"
$ws.Range("E2").Value = "'1.37%"
$ws.Range("D3").Value = "'29.61"
$ws.Range("E3").Value = "'4.22%"
$ws.Range("D4").Value = "'5.130"
$ws.Range("E4").Value = "'1.50%"
$ws.Range("E5").Value = "'3.24%"
$ws.Range("D6").Value = "'7.335"
$ws.Range("E6").Value = "'1.54%"
$ws.Range("D7").Value = "'3.395"
$ws.Range("E7").Value = "'0.89%"
$ws.Range("D8").Value = "'1.367"
$ws.Range("E8").Value = "'-2.37%"
$ws.Range("D9").Value = "'0.9206"
$ws.Range("E9").Value = "'0.33%"
$ws.Range("D10").Value = "'0.1592"
$ws.Range("E10").Value = "'3.56%"
$ws.Range("D11").Value = "'0.06815"
$ws.Range("E11").Value = "'2.38%"
$ws.Range("D12").Value = "'0.07707"
$ws.Range("E12").Value = "'1.23%"
$ws.Range("D13").Value = "'0.02936"
$ws.Range("E13").Value = "'5.51%"
$ws.Range("D14").Value = "'0.08985"
$ws.Range("E14").Value = "'0.16%"
$ws.Range("D15").Value = "'0.001591"
$ws.Range("E15").Value = "'0.07%"
$ws.Range("D16").Value = "'0.04481"
$ws.Range("E16").Value = "'1.15%"
$ws.Range("D17").Value = "'0.0006453"
$ws.Range("E17").Value = "'1.71%"
$ws.Range("D18").Value = "'0.006248"
$ws.Range("E18").Value = "'1.59%"
$ws.Range("E19").Value = "'-0.16%"
$ws.Range("D20").Value = "'2.228"
$ws.Range("E20").Value = "'-0.61%"
$ws.Range("D21").Value = "'0.3215"
$ws.Range("E21").Value = "'1.04%"
$ws.Range("D22").Value = "'0.1303"
$ws.Range("E22").Value = "'-2.87%"
$ws.Range("D23").Value = "'4.062"
$ws.Range("E23").Value = "'1.15%"
$ws.Range("E24").Value = "'2.40%"
$ws.Range("D25").Value = "'0.001191"
$ws.Range("E25").Value = "'0.74%"
$ws.Range("D26").Value = "'0.004121"
$ws.Range("E26").Value = "'-7.59%"
$ws.Range("E27").Value = "'-0.05%"
$ws.Range("D28").Value = "'0.0001615"
$ws.Range("E28").Value = "'-0.20%"
$ws.Range("D40").Value = "'0.04274"
$ws.Range("E40").Value = "'3.86%"
$ws.Range("D41").Value = "'0.006709"
$ws.Range("E41").Value = "'1.44%"
$ws.Range("E42").Value = "'0.54%"
$ws.Range("D43").Value = "'0.002197"
$ws.Range("E43").Value = "'6.74%"
$ws.Range("E44").Value = "'3.71%"
$ws.Range("D45").Value = "'0.00005670"
$ws.Range("E45").Value = "'5.09%"
$ws.Range("D46").Value = "'1.974"
$ws.Range("E46").Value = "'2.14%"
$ws.Range("E47").Value = "'-29.50%"
